$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D and E contain numeric-looking text (prices and
# percentages) that must remain plain text, exactly as scraped from the
# source site. Force the Text number format before assigning so Excel
# does not coerce the values into numbers.
$textCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6",
    "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11",
    "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16",
    "D17", "E17", "D18", "E18", "E19", "D20", "E20", "E21", "D22", "E22",
    "D23", "E23", "D24", "E24", "D25", "E25", "E26", "D27", "D39", "E39",
    "D40", "E40", "D41", "E41", "E42", "D43", "E43", "E44", "D45", "E45",
    "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "E50", "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, matching the source diff.
$ws.Range("D2").Value = "311.54"
$ws.Range("E2").Value = "-0.47%"
$ws.Range("D3").Value = "37.64"
$ws.Range("E3").Value = "-1.31%"
$ws.Range("D4").Value = "5.167"
$ws.Range("E4").Value = "1.55%"
$ws.Range("D5").Value = "0.07926"
$ws.Range("E5").Value = "2.17%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "4.430"
$ws.Range("E6").Value = "1.75%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "1.933"
$ws.Range("E7").Value = "0.95%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "8.281"
$ws.Range("E8").Value = "1.21%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "2.989"
$ws.Range("E9").Value = "-3.54%"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "0.9286"
$ws.Range("E10").Value = "1.33%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "0.1126"
$ws.Range("E11").Value = "-8.82%"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "0.1921"
$ws.Range("E12").Value = "1.46%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "0.09078"
$ws.Range("E13").Value = "2.07%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03300"
$ws.Range("E14").Value = "-2.58%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09610"
$ws.Range("E15").Value = "-1.04%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001380"
$ws.Range("E16").Value = "0.43%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "0.005866"
$ws.Range("E17").Value = "-3.16%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.596"
$ws.Range("E18").Value = "1.69%"
$ws.Range("E19").Value = "0.06%"
$ws.Range("D20").Value = "5.951"
$ws.Range("E20").Value = "18.31%"
$ws.Range("E21").Value = "0.60%"
$ws.Range("D22").Value = "0.2590"
$ws.Range("E22").Value = "-0.03%"
$ws.Range("D23").Value = "0.04414"
$ws.Range("E23").Value = "0.31%"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").Value = "1.78%"
$ws.Range("D25").Value = "0.004633"
$ws.Range("E25").Value = "9.07%"
$ws.Range("E26").Value = "0.72%"
$ws.Range("D27").Value = "0.0003991"
$ws.Range("D39").Value = "0.02252"
$ws.Range("E39").Value = "5.64%"
$ws.Range("D40").Value = "0.05094"
$ws.Range("E40").Value = "1.80%"
$ws.Range("D41").Value = "0.007477"
$ws.Range("E41").Value = "-4.88%"
$ws.Range("E42").Value = "-10.21%"
$ws.Range("D43").Value = "0.1356"
$ws.Range("E43").Value = "0.76%"
$ws.Range("E44").Value = "3.38%"
$ws.Range("D45").Value = "0.008627"
$ws.Range("E45").Value = "-10.93%"
$ws.Range("D46").Value = "0.00006633"
$ws.Range("E46").Value = "1.57%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "-0.03%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "0.001000"
$ws.Range("E48").Value = "-40.76%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "0.002870"
$ws.Range("E49").Value = "-10.31%"
$ws.Range("E50").Value = "-0.03%"
$ws.Range("E51").Value = "-0.03%"

Write-Output "Applied cryptos.xlsx symbol/price update"